# Daily attendance processing - 2025-12-13 05:53:55
# Reverses the order of comma-separated "Recorded By" entries in column G
# for every data row on the active sheet (only affects cells that actually
# contain more than one comma-separated value; single-value cells are
# left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversedParts = $parts[($parts.Length - 1)..0]
        $cell.Value2 = [string]::Join(", ", $reversedParts)
    }
}
